$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
  @(2, "TBA323121451756", "BOG0011142986", "2025-08-08 03:57:05", "NO ESCANEADO", "2025-08-08 04:19:21"),
  @(3, "TBA323119682119", "BOG0011142997", "2025-08-08 03:57:05", "NO ESCANEADO", "2025-08-08 04:10:26"),
  @(4, "TBA323098427821", "BOG0011142963", "2025-08-08 03:57:05", "NO ESCANEADO", "2025-08-08 04:10:26"),
  @(5, "TBA323098408870", "BOG0011142965", "2025-08-08 03:57:05", "NO ESCANEADO", "2025-08-08 04:10:26"),
  @(6, "TBA323098354329", "BOG0011142951", "2025-08-08 03:57:05", "NO ESCANEADO", "2025-08-08 04:10:26"),
  @(7, "TBA323098196734", "BOG0011142955", "2025-08-08 03:57:05", "NO ESCANEADO", "2025-08-08 04:10:26"),
  @(8, "TBA323098175547", "BOG0011142962", "2025-08-08 03:57:05", "NO ESCANEADO", "2025-08-08 04:10:26"),
  @(9, "TBA323096312137", "BOG0011143002", "2025-08-08 03:57:05", "NO ESCANEADO", "2025-08-08 04:10:26"),
  @(10, "TBA323096042734", "BOG0011142912", "2025-08-08 03:57:05", "NO ESCANEADO", "2025-08-08 04:10:26"),
  @(11, "TBA323094564948", "BOG0011142893", "2025-08-08 03:57:05", "NO ESCANEADO", "2025-08-08 04:10:26"),
  @(12, "TBA323088149841", "BOG0011142930", "2025-08-08 03:57:05", "NO ESCANEADO", "2025-08-08 04:10:26"),
  @(13, "TBA323087094055", "BOG0011142940", "2025-08-08 03:57:05", "NO ESCANEADO", "2025-08-08 04:10:26"),
  @(14, "TBA323086788090", "BOG0011142957", "2025-08-08 03:57:05", "NO ESCANEADO", "2025-08-08 04:10:26"),
  @(15, "TBA323086407252", "BOG0011142954", "2025-08-08 03:57:05", "NO ESCANEADO", "2025-08-08 04:10:26"),
  @(16, "TBA323086147012", "BOG0011142959", "2025-08-08 03:57:05", "NO ESCANEADO", "2025-08-08 04:10:26"),
  @(17, "TBA323086107431", "BOG0011142886", "2025-08-08 03:57:05", "NO ESCANEADO", "2025-08-08 04:10:26"),
  @(18, "TBA323086047239", "BOG0011142926", "2025-08-08 03:57:05", "NO ESCANEADO", "2025-08-08 04:10:26"),
  @(19, "TBA323085896044", "BOG0011142968", "2025-08-08 03:57:05", "NO ESCANEADO", "2025-08-08 04:10:26"),
  @(20, "TBA323085574582", "BOG0011142960", "2025-08-08 03:57:05", "NO ESCANEADO", "2025-08-08 04:10:26"),
  @(21, "TBA323085552381", "BOG0011142905", "2025-08-08 03:57:05", "NO ESCANEADO", "2025-08-08 04:10:26"),
  @(22, "TBA323079856372", "BOG0011142843", "2025-08-08 03:57:05", "NO ESCANEADO", "2025-08-08 04:10:26"),
  @(23, "TBA323079830678", "BOG0011142939", "2025-08-08 03:57:05", "NO ESCANEADO", "2025-08-08 04:10:26"),
  @(24, "TBA323079192344", "BOG0011142845", "2025-08-08 03:57:05", "NO ESCANEADO", "2025-08-08 04:10:26"),
  @(25, "TBA323079103133", "BOG0011142952", "2025-08-08 03:57:05", "NO ESCANEADO", "2025-08-08 04:10:26"),
  @(26, "TBA323078858716", "BOG0011142934", "2025-08-08 03:57:05", "NO ESCANEADO", "2025-08-08 04:10:26"),
  @(27, "TBA323078134959", "BOG0011142933", "2025-08-08 03:57:05", "NO ESCANEADO", "2025-08-08 04:10:26"),
  @(28, "TBA323076853839", "BOG0011142937", "2025-08-08 03:57:05", "NO ESCANEADO", "2025-08-08 04:10:26"),
  @(29, "TBA323074861015", "BOG0011142953", "2025-08-08 03:57:05", "NO ESCANEADO", "2025-08-08 04:10:26"),
  @(30, "TBA323073407805", "BOG0011142938", "2025-08-08 03:57:05", "NO ESCANEADO", "2025-08-08 04:10:26"),
  @(31, "TBA323073116420", "BOG0011142976", "2025-08-08 03:57:05", "NO ESCANEADO", "2025-08-08 04:10:26"),
  @(32, "TBA323071837440", "BOG0011142966", "2025-08-08 03:57:05", "NO ESCANEADO", "2025-08-08 04:10:26"),
  @(33, "TBA323071578041", "BOG0011142942", "2025-08-08 03:57:05", "NO ESCANEADO", "2025-08-08 04:10:26"),
  @(34, "TBA323071019079", "BOG0011142802", "2025-08-08 03:57:05", "NO ESCANEADO", "2025-08-08 04:10:26"),
  @(35, "TBA323070013397", "BOG0011142800", "2025-08-08 03:57:05", "NO ESCANEADO", "2025-08-08 04:10:26"),
  @(36, "TBA323068851911", "BOG0011142935", "2025-08-08 03:57:05", "NO ESCANEADO", "2025-08-08 04:10:26"),
  @(37, "TBA323068796349", "BOG0011142795", "2025-08-08 03:57:05", "NO ESCANEADO", "2025-08-08 04:10:26"),
  @(38, "TBA323068278341", "BOG0011142775", "2025-08-08 03:57:05", "NO ESCANEADO", "2025-08-08 04:10:26"),
  @(39, "TBA323068014565", "BOG0011142878", "2025-08-08 03:57:05", "NO ESCANEADO", "2025-08-08 04:10:26"),
  @(40, "TBA323067483217", "BOG0011142791", "2025-08-08 03:57:05", "NO ESCANEADO", "2025-08-08 04:10:26"),
  @(41, "TBA323067481425", "BOG0011142731", "2025-08-08 03:57:05", "NO ESCANEADO", "2025-08-08 04:10:26"),
  @(42, "TBA323067379935", "BOG0011142733", "2025-08-08 03:57:05", "NO ESCANEADO", "2025-08-08 04:10:26"),
  @(43, "TBA323067215546", "BOG0011142724", "2025-08-08 03:57:05", "NO ESCANEADO", "2025-08-08 04:10:26"),
  @(44, "TBA323066955577", "BOG0011142932", "2025-08-08 03:57:05", "NO ESCANEADO", "2025-08-08 04:10:26"),
  @(45, "TBA323065748432", "BOG0011142943", "2025-08-08 03:57:05", "NO ESCANEADO", "2025-08-08 04:10:26"),
  @(46, "TBA323065514270", "BOG0011142769", "2025-08-08 03:57:05", "NO ESCANEADO", "2025-08-08 04:10:26"),
  @(47, "TBA323065354447", "BOG0011142788", "2025-08-08 03:57:05", "NO ESCANEADO", "2025-08-08 04:10:26"),
  @(48, "TBA323065189260", "BOG0011142716", "2025-08-08 03:57:05", "NO ESCANEADO", "2025-08-08 04:10:26"),
  @(49, "TBA323064949443", "BOG0011142740", "2025-08-08 03:57:05", "NO ESCANEADO", "2025-08-08 04:10:26"),
  @(50, "TBA323064833422", "BOG0011142721", "2025-08-08 03:57:05", "NO ESCANEADO", "2025-08-08 04:10:26"),
  @(51, "TBA323064773485", "BOG0011142785", "2025-08-08 03:57:05", "NO ESCANEADO", "2025-08-08 04:10:26"),
  @(52, "TBA323064442166", "BOG0011142794", "2025-08-08 03:57:05", "NO ESCANEADO", "2025-08-08 04:10:26"),
  @(53, "TBA323060264752", "BOG0011142778", "2025-08-08 03:57:05", "NO ESCANEADO", "2025-08-08 04:10:26"),
  @(54, "TBA323060200298", "BOG0011142782", "2025-08-08 03:57:05", "NO ESCANEADO", "2025-08-08 04:10:26"),
  @(55, "TBA323060175410", "BOG0011142803", "2025-08-08 03:57:05", "NO ESCANEADO", "2025-08-08 04:10:26"),
  @(56, "TBA323060080360", "BOG0011142936", "2025-08-08 03:57:05", "NO ESCANEADO", "2025-08-08 04:10:26"),
  @(57, "TBA323059932389", "BOG0011142735", "2025-08-08 03:57:05", "NO ESCANEADO", "2025-08-08 04:10:26"),
  @(58, "TBA323056168803", "BOG0011142790", "2025-08-08 03:57:05", "NO ESCANEADO", "2025-08-08 04:10:26"),
  @(59, "TBA323054173890", "BOG0011142797", "2025-08-08 03:57:05", "NO ESCANEADO", "2025-08-08 04:10:26"),
  @(60, "TBA323054042345", "BOG0011142792", "2025-08-08 03:57:05", "NO ESCANEADO", "2025-08-08 04:10:26"),
  @(61, "TBA323051817440", "BOG0011142779", "2025-08-08 03:57:05", "NO ESCANEADO", "2025-08-08 04:10:26"),
  @(62, "TBA323050760089", "BOG0011142804", "2025-08-08 03:57:05", "NO ESCANEADO", "2025-08-08 04:10:26"),
  @(63, "TBA323050436599", "BOG0011142548", "2025-08-08 03:57:05", "NO ESCANEADO", "2025-08-08 04:10:26"),
  @(64, "TBA323050306127", "BOG0011142635", "2025-08-08 03:57:05", "NO ESCANEADO", "2025-08-08 04:10:26"),
  @(65, "TBA323049275611", "BOG0011142801", "2025-08-08 03:57:05", "NO ESCANEADO", "2025-08-08 04:10:26"),
  @(66, "TBA323048756559", "BOG0011142789", "2025-08-08 03:57:05", "NO ESCANEADO", "2025-08-08 04:10:26"),
  @(67, "TBA323047725294", "BOG0011142627", "2025-08-08 03:57:05", "NO ESCANEADO", "2025-08-08 04:10:26"),
  @(68, "TBA323046640835", "BOG0011142781", "2025-08-08 03:57:05", "NO ESCANEADO", "2025-08-08 04:10:26"),
  @(69, "TBA323046489893", "BOG0011142639", "2025-08-08 03:57:05", "NO ESCANEADO", "2025-08-08 04:10:26"),
  @(70, "TBA323045835654", "BOG0011142614", "2025-08-08 03:57:05", "NO ESCANEADO", "2025-08-08 04:10:26"),
  @(71, "TBA323045525229", "BOG0011142592", "2025-08-08 03:57:05", "NO ESCANEADO", "2025-08-08 04:10:26"),
  @(72, "TBA323045331849", "BOG0011142640", "2025-08-08 03:57:05", "NO ESCANEADO", "2025-08-08 04:10:26"),
  @(73, "TBA323044998566", "BOG0011142642", "2025-08-08 03:57:05", "NO ESCANEADO", "2025-08-08 04:10:26"),
  @(74, "TBA323044603885", "BOG0011142616", "2025-08-08 03:57:05", "NO ESCANEADO", "2025-08-08 04:10:26"),
  @(75, "TBA323044225455", "BOG0011142624", "2025-08-08 03:57:05", "NO ESCANEADO", "2025-08-08 04:10:26"),
  @(76, "TBA323043605568", "BOG0011142605", "2025-08-08 03:57:05", "NO ESCANEADO", "2025-08-08 04:10:26"),
  @(77, "TBA323043396001", "BOG0011142641", "2025-08-08 03:57:05", "NO ESCANEADO", "2025-08-08 04:10:26"),
  @(78, "TBA323042863218", "BOG0011142608", "2025-08-08 03:57:05", "NO ESCANEADO", "2025-08-08 04:10:26"),
  @(79, "TBA323042358097", "BOG0011142771", "2025-08-08 03:57:05", "NO ESCANEADO", "2025-08-08 04:10:26"),
  @(80, "TBA323041079098", "BOG0011142591", "2025-08-08 03:57:05", "NO ESCANEADO", "2025-08-08 04:10:26"),
  @(81, "TBA323040906555", "BOG0011142622", "2025-08-08 03:57:05", "NO ESCANEADO", "2025-08-08 04:10:26"),
  @(82, "TBA323040840858", "BOG0011142770", "2025-08-08 03:57:05", "NO ESCANEADO", "2025-08-08 04:10:26"),
  @(83, "TBA323040511974", "BOG0011142765", "2025-08-08 03:57:05", "NO ESCANEADO", "2025-08-08 04:10:26"),
  @(84, "TBA323040436706", "BOG0011142643", "2025-08-08 03:57:05", "NO ESCANEADO", "2025-08-08 04:10:26"),
  @(85, "TBA323040228272", "BOG0011142621", "2025-08-08 03:57:05", "NO ESCANEADO", "2025-08-08 04:10:26"),
  @(86, "TBA323039927821", "BOG0011142663", "2025-08-08 03:57:05", "NO ESCANEADO", "2025-08-08 04:10:26"),
  @(87, "TBA323039720003", "BOG0011142571", "2025-08-08 03:57:05", "NO ESCANEADO", "2025-08-08 04:10:26"),
  @(88, "TBA323039687648", "BOG0011142623", "2025-08-08 03:57:05", "NO ESCANEADO", "2025-08-08 04:10:26"),
  @(89, "TBA323039654387", "BOG0011142777", "2025-08-08 03:57:05", "NO ESCANEADO", "2025-08-08 04:10:26"),
  @(90, "TBA323039441832", "BOG0011142787", "2025-08-08 03:57:05", "NO ESCANEADO", "2025-08-08 04:10:26"),
  @(91, "TBA323039298468", "BOG0011142617", "2025-08-08 03:57:05", "NO ESCANEADO", "2025-08-08 04:10:26"),
  @(92, "TBA323039084420", "BOG0011142625", "2025-08-08 03:57:05", "NO ESCANEADO", "2025-08-08 04:10:26"),
  @(93, "TBA323039003113", "BOG0011142647", "2025-08-08 03:57:05", "NO ESCANEADO", "2025-08-08 04:10:26"),
  @(94, "TBA323038796730", "BOG0011142633", "2025-08-08 03:57:05", "NO ESCANEADO", "2025-08-08 04:10:26"),
  @(95, "TBA323038688439", "BOG0011142648", "2025-08-08 03:57:05", "NO ESCANEADO", "2025-08-08 04:10:26"),
  @(96, "TBA323038429030", "BOG0011142618", "2025-08-08 03:57:05", "NO ESCANEADO", "2025-08-08 04:10:26"),
  @(97, "TBA323037911002", "BOG0011142645", "2025-08-08 03:57:05", "NO ESCANEADO", "2025-08-08 04:10:26"),
  @(98, "TBA323067752726", "BOG0011142784", "2025-08-08 03:57:05", "RECIBIDO", "2025-08-08 04:08:12"),
  @(99, "TBA323046215869", "BOG0011142646", "2025-08-08 03:57:05", "RECIBIDO", "2025-08-08 04:05:04"),
  @(100, "TBA323047693856", "BOG0011142737", "2025-08-08 03:57:05", "RECIBIDO", "2025-08-08 04:04:54"),
  @(101, "789", "N/A", "2025-08-08 04:04:25", "NO ESPERADO", "2025-08-08 04:04:25"),
  @(102, "4657", "N/A", "2025-08-08 04:02:42", "NO ESPERADO", "2025-08-08 04:02:42"),
  @(103, "14567", "N/A", "2025-08-08 04:00:19", "NO ESPERADO", "2025-08-08 04:00:19"),
  @(104, "1234", "N/A", "2025-08-08 03:59:54", "NO ESPERADO", "2025-08-08 03:59:54"),
  @(105, "TBA323043552239", "BOG0011142691", "2025-08-08 03:57:05", "RECIBIDO", "2025-08-08 03:59:46"),
  @(106, "TBA323041521274", "BOG0011142780", "2025-08-08 03:57:05", "RECIBIDO", "2025-08-08 03:59:17"),
  @(107, "TBA323039982955", "BOG0011142607", "2025-08-08 03:57:05", "RECIBIDO", "2025-08-08 03:59:10"),
  @(108, "TBA323039355408", "BOG0011142774", "2025-08-08 03:57:05", "RECIBIDO", "2025-08-08 03:59:01"),
  @(109, "TBA323038545754", "BOG0011142620", "2025-08-08 03:57:05", "RECIBIDO", "2025-08-08 03:58:53")
)

$textForceRows = @(101, 102, 103, 104)

foreach ($row in $rows) {
    $r = $row[0]
    if ($textForceRows -contains $r) {
        $ws.Cells.Item($r, 1).NumberFormat = "@"
        $ws.Cells.Item($r, 1).Value = $row[1]
        $ws.Cells.Item($r, 1).Style = "Normal"
    } else {
        $ws.Cells.Item($r, 1).Value = $row[1]
    }
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
    $ws.Cells.Item($r, 4).Value = $row[4]
    $ws.Cells.Item($r, 5).Value = $row[5]
}
